$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.467.66"
$ws.Range("E2").Value = "  +6.13%  "
$ws.Range("D3").Value = "2.475.73"
$ws.Range("E3").Value = "  +7.60%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'570.98"
$ws.Range("E5").Value = "  +5.57%  "
$ws.Range("D6").Value = "'143.80"
$ws.Range("E6").Value = "  +11.87%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("E8").Value = "  +4.32%  "
$ws.Range("D9").Value = "2.474.82"
$ws.Range("E9").Value = "  +7.66%  "
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("D11").Value = "'5.77"
$ws.Range("E11").Value = "  +4.57%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D14").Value = "'26.48"
$ws.Range("E14").Value = "  +14.80%  "
$ws.Range("D15").Value = "2.914.25"
$ws.Range("E15").Value = "  +7.34%  "
$ws.Range("D16").Value = "63.328.10"
$ws.Range("E16").Value = "  +6.01%  "
$ws.Range("E17").Value = "  +10.66%  "
$ws.Range("D18").Value = "2.471.89"
$ws.Range("E18").Value = "  +7.37%  "
$ws.Range("D19").Value = "'11.36"
$ws.Range("E19").Value = "  +9.38%  "
$ws.Range("D20").Value = "'345.49"
$ws.Range("E20").Value = "  +11.56%  "
$ws.Range("D21").Value = "'4.34"
$ws.Range("E21").Value = "  +8.27%  "
$ws.Range("D22").Value = "'6.86"
$ws.Range("E22").Value = "  +6.47%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "'65.84"
$ws.Range("E24").Value = "  +4.11%  "
$ws.Range("E25").Value = "  +4.13%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +13.85%  "
$ws.Range("E28").Value = "  +7.03%  "
$ws.Range("D29").Value = "'1.31"
$ws.Range("E29").Value = "  +11.33%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0818"
$ws.Range("E30").Value = "  +15.29%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.85"
$ws.Range("E31").Value = "  +18.55%  "
$ws.Range("D32").Value = "'1.85"
$ws.Range("E32").Value = "  +8.90%  "
$ws.Range("D33").Value = "'175.11"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").Value = "'1.51"
$ws.Range("E34").Value = "  +12.85%  "
$ws.Range("E35").Value = "  +6.36%  "
$ws.Range("D36").Value = "'18.95"
$ws.Range("E36").Value = "  +7.30%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'372.09"
$ws.Range("E37").Value = "  +19.98%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'4.51"
$ws.Range("E38").Value = "  +11.52%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'1.72"
$ws.Range("E41").Value = "  +15.39%  "
$ws.Range("D42").Value = "'40.31"
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("D43").Value = "'151.77"
$ws.Range("E43").Value = "  +11.78%  "
$ws.Range("E44").Value = "  +10.17%  "
$ws.Range("E45").Value = "  +12.84%  "
$ws.Range("D46").Value = "'0.602"
$ws.Range("E46").Value = "  +6.71%  "
$ws.Range("E48").Value = "  +7.84%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0238"
$ws.Range("E49").Value = "  +6.52%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0227"
$ws.Range("E50").Value = "  +7.46%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'18.34"
$ws.Range("E51").Value = "  +10.70%  "
